$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Bump the published Version (row 3) and regenerate the Date (row 8).
$ws.Range("B3").Value = "0.2.0"
$ws.Range("B8").Value = "2023-10-20T08:59:58+00:00"

# A new "Jurisdiction" property needs to land right before "Description", which pushes
# Description / Purpose / Copyright / Immutable down by one row (rows 11-14 -> 12-15).
# Capture the old row contents first (read via the indexer so we get real values, not
# property descriptors), then write them back one row lower, bottom row first so nothing
# gets clobbered before it's been read.
$a11 = $ws.Range("A11").Value()
$b11 = $ws.Range("B11").Value()
$a12 = $ws.Range("A12").Value()
$b12 = $ws.Range("B12").Value()
$a13 = $ws.Range("A13").Value()
$b13 = $ws.Range("B13").Value()
$a14 = $ws.Range("A14").Value()
$b14 = $ws.Range("B14").Value()

# The new row 15 needs the same look as the rest of the table, so clone row 14's format
# into it before writing values there.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

$ws.Range("A15").Value = $a14
$ws.Range("B15").Value = $b14
$ws.Range("A14").Value = $a13
$ws.Range("B14").Value = $b13
$ws.Range("A13").Value = $a12
$ws.Range("B13").Value = $b12
$ws.Range("A12").Value = $a11
$ws.Range("B12").Value = $b11

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
